$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for the two new columns
$ws.Range("F1").Value = "actual_conductivity"
$ws.Range("G1").Value = "temp"

# actual_conductivity values (F2:F6)
$ws.Range("F2").Value = 26.57
$ws.Range("F3").Value = 37.09
$ws.Range("F4").Value = 484.9
$ws.Range("F5").Value = 58.86
$ws.Range("F6").Value = 57.59

# temp values (G2:G6)
$ws.Range("G2").Value = 14
$ws.Range("G3").Value = 12.9
$ws.Range("G4").Value = 13.2
$ws.Range("G5").Value = 13.4
$ws.Range("G6").Value = 12.7

# Copy the formatting of an existing styled header cell onto the new headers
# so they reuse the existing style index instead of creating a new one.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("F1:G1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Widen column F to fit the new header text
$ws.Columns.Item(6).ColumnWidth = 16.14

# Restore the selection like the source workbook
$ws.Range("J8").Select() | Out-Null
